$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "243.65"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.88%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "27.03"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "4.31%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.153"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.72%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05619"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.53%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.489"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.10%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8169"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.09%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8320"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.91%"

$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1329"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.81%"

$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06901"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-0.93%"

$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "BitrueCoin"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.02891"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.38%"

$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "BitMartToken"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09374"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.29%"

$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "BitForexToken"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.001511"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.22%"

$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "One"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0005989"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-93.83%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.006115"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-2.45%"

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.64%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.022"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.21%"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.307"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "8.91%"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.71%"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.03092"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-3.93%"

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-2.18%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.741"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "0.06%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04559"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-3.08%"

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-2.50%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001225"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-2.01%"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004487"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-2.55%"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.00009797"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "2.05%"

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "0.65%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03635"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-0.51%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006068"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.85%"

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-0.26%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002590"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "4.32%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008178"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "3.78%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005310"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "0.00%"

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.03%"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-18.37%"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002655"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "29.67%"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002099"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.03%"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0001999"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.03%"
